$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force text storage for the data cells (so numeric-looking / percent-looking
# strings like "23" or "0%" keep being stored as shared-string text, matching
# the original workbook's cell typing) and then restore the original
# (default/general) number format so the cell style index is unchanged.
$dataRange = $ws.Range("A2:F6")
$dataRange.NumberFormat = "@"

$ws.Range("A2").Value = "02/03/2023 12:12:39"
$ws.Range("B2").Value = "Kansas City"
$ws.Range("C2").Value = "23"
$ws.Range("D2").Value = "0%"
$ws.Range("E2").Value = "7 mph"
$ws.Range("F2").Value = "13 mph"

$ws.Range("A3").Value = "02/03/2023 12:12:42"
$ws.Range("B3").Value = "New York"
$ws.Range("C3").Value = "24"
$ws.Range("D3").Value = "0%"
$ws.Range("E3").Value = "13 mph"
$ws.Range("F3").Value = "13 mph"

$ws.Range("A4").Value = "02/03/2023 12:12:44"
$ws.Range("B4").Value = "Sacramento"
$ws.Range("C4").Value = "46"
$ws.Range("D4").Value = "77%"
$ws.Range("E4").Value = "7 mph"
$ws.Range("F4").Value = "13 mph"

$ws.Range("A5").Value = "02/03/2023 12:12:47"
$ws.Range("B5").Value = "Chicago"
$ws.Range("C5").Value = "8"
$ws.Range("D5").Value = "1%"
$ws.Range("E5").Value = "7 mph"
$ws.Range("F5").Value = "13 mph"

$ws.Range("A6").Value = "02/03/2023 12:12:50"
$ws.Range("B6").Value = "Nashville"
$ws.Range("C6").Value = "31"
$ws.Range("D6").Value = "0%"
$ws.Range("E6").Value = "9 mph"
$ws.Range("F6").Value = "13 mph"

$dataRange.NumberFormat = ""
